# Update "社会融资规模及构成" sheet: refresh the data table with 2010-2022
# figures (replacing the old 2002-2020 figures) and drop the now-unused
# trailing rows so the sheet dimension shrinks from A1:I20 to A1:I14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=year, B=人民币贷款, C=企业债券, D=信托贷款, E=外币贷款(折合人民币),
#          F=委托贷款, G=未贴现银行承兑汇票, H=社会融资规模, I=非金融企业境内股票
$rows = @(
  ,@("2010年", 79451,  11063, 3865,   4855,   8748,   23346,   140191, 5786)
  ,@("2011年", 74715,  13658, 2034,   5712,   12962,  10271,   128286, 4377)
  ,@("2012年", 82038,  22551, 12845,  9163,   12838,  10499,   157631, 2508)
  ,@("2013年", 88916,  18111, 18404,  5848,   25466,  7756,    173169, 2219)
  ,@("2014年", 97452,  24329, 5174,   1235,   21740,  -1198,   158761, 4350)
  ,@("2015年", 112693, 29388, 434,    -6427,  15911,  -10567,  154063, 7590)
  ,@("2016年", 124372, 29865, 8593,   -5640,  21854,  -19514,  177999, 12416)
  ,@("2017年", 138432, 6244,  22232,  18,     7994,   5364,    261536, 8759)
  ,@("2018年", 156712, 26318, -6975,  -4201,  -16062, -6343,   224920, 3606)
  ,@("2019年", 168835, 33384, -3467,  -1275,  -9396,  -4757,   256735, 3479)
  ,@("2020年", 200310, 43748, -11020, 1450,   -3954,  1746,    347917, 8923)
  ,@("2021年", 199403, 32866, -20074, 1715,   -1696,  -4916,   313407, 12133)
  ,@("2022年", 209149, 20508, -6003,  -5254,  3579,   -3411,   320101, 11757)
)

$r = 2
foreach ($row in $rows) {
  for ($c = 0; $c -lt $row.Length; $c++) {
    $ws.Cells.Item($r, $c + 1).Value = $row[$c]
  }
  $r++
}

# The new table ends at row 14 (2022年); the old sheet had data through row
# 20 (2020年) so the trailing rows 15:20 are no longer needed.
$ws.Rows("15:20").Delete()
